$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 10 new workout rows (270-279) below the existing data (which ends
# at row 269), reusing the date-formatted style already used by column B
# in the row above so no new cell style is introduced.
# ---------------------------------------------------------------------------

$lastExistingRow = 269
$newRowCount = 10

# Copy the formatting of the last existing data row down onto the new rows
# so that column B keeps its date number format (style index 1) and the
# other columns keep the default style - matches the source file exactly.
$srcFormatRange = $ws.Range("A" + $lastExistingRow + ":M" + $lastExistingRow)
$destFormatRange = $ws.Range("A270:M279")
$srcFormatRange.Copy()
$destFormatRange.PasteSpecial(-4122)

# Columns: A Participant, B Date, C Workout Type, D Total Duration,
#          E Total Distance, F Total Elevation, G Zone1, H Zone2, I Zone3,
#          J Zone4, K Zone5, L Workout Level, M Week
$rows = @(
    @{ A="Jeremiah"; B=45497; C="Run";     D=23; E=2.46; F=133; G=0;  H=16; I=5;  J=0; K=0; L="Sauntering Hippo"; M=7 },
    @{ A="Steven";   B=45497; C="Workout"; D=64; E=0;    F=0;   G=44; H=19; I=2;  J=0; K=0; L="Brave Leopard";    M=7 },
    @{ A="Steven";   B=45497; C="Walk";    D=33; E=1.73; F=82;  G=33; H=0;  I=0;  J=0; K=0; L="Brave Leopard";    M=7 },
    @{ A="Steven";   B=45497; C="Walk";    D=45; E=1.9;  F=66;  G=45; H=0;  I=0;  J=0; K=0; L="Brave Leopard";    M=7 },
    @{ A="Steven";   B=45497; C="Walk";    D=29; E=1.45; F=69;  G=29; H=0;  I=0;  J=0; K=0; L="Brave Leopard";    M=7 },
    @{ A="Matt";     B=45497; C="Walk";    D=25; E=0.93; F=69;  G=25; H=0;  I=0;  J=0; K=0; L="Agile Antelope";   M=7 },
    @{ A="Phil";     B=45497; C="Workout"; D=79; E=0;    F=0;   G=55; H=23; I=2;  J=0; K=0; L="Sauntering Hippo"; M=7 },
    @{ A="Steven";   B=45498; C="Run";     D=36; E=3.05; F=69;  G=2;  H=7;  I=24; J=1; K=0; L="Brave Leopard";    M=7 },
    @{ A="Steven";   B=45498; C="Walk";    D=17; E=0.85; F=49;  G=17; H=0;  I=0;  J=0; K=0; L="Brave Leopard";    M=7 },
    @{ A="Steven";   B=45498; C="Walk";    D=20; E=1.09; F=56;  G=20; H=0;  I=0;  J=0; K=0; L="Brave Leopard";    M=7 }
)

$r = $lastExistingRow
foreach ($row in $rows) {
    $r = $r + 1
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
}

# ---------------------------------------------------------------------------
# Update the sheet view: keep the header row frozen, scroll further down so
# row 255 is the first visible row below the freeze, and leave the final
# selection on A280 (the first empty row after the new data).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A255").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 255

$ws.Range("A280").Select()
